$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 406
$ws.Range("F3").Value = 1087
$ws.Range("F4").Value = 9663
$ws.Range("F5").Value = 208
$ws.Range("F8").Value = 6587
$ws.Range("F9").Value = 632
$ws.Range("F10").Value = 10487
$ws.Range("F11").Value = 11645
$ws.Range("F12").Value = 1259
$ws.Range("F13").Value = 1204
$ws.Range("F14").Value = 5049
$ws.Range("F15").Value = 837
$ws.Range("F16").Value = 502
$ws.Range("F20").Value = 1371
$ws.Range("F21").Value = 282
$ws.Range("F22").Value = 1911
$ws.Range("F23").Value = 924
$ws.Range("F24").Value = 1327
$ws.Range("F25").Value = 864
$ws.Range("F26").Value = 7
$ws.Range("F27").Value = 2081
$ws.Range("G27").Value = 80
$ws.Range("F28").Value = 445
$ws.Range("F29").Value = 667
$ws.Range("F30").Value = 2757
$ws.Range("F32").Value = 1840
$ws.Range("F34").Value = 837
$ws.Range("F35").Value = 90
$ws.Range("F36").Value = 938
$ws.Range("F37").Value = 31
$ws.Range("F38").Value = 56
$ws.Range("F39").Value = 3438
$ws.Range("F41").Value = 94
$ws.Range("F42").Value = 533
$ws.Range("F43").Value = 599
$ws.Range("F45").Value = 904
$ws.Range("F47").Value = 11
$ws.Range("F48").Value = 4239
$ws.Range("F49").Value = 88

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 10
$ws.Range("F9").Value = 33
$ws.Range("F25").Value = 107

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 6115

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 406
$ws.Range("F3").Value = 1087
$ws.Range("F4").Value = 9663
$ws.Range("F7").Value = 632
$ws.Range("F9").Value = 10487
$ws.Range("F10").Value = 11646
$ws.Range("F12").Value = 1204
$ws.Range("F13").Value = 5049
$ws.Range("F14").Value = 837
$ws.Range("F15").Value = 502
$ws.Range("F18").Value = 33
$ws.Range("F20").Value = 1371
$ws.Range("F21").Value = 282
$ws.Range("F22").Value = 1911
$ws.Range("F23").Value = 924
$ws.Range("F24").Value = 1327
$ws.Range("F25").Value = 864
$ws.Range("F26").Value = 2081
$ws.Range("G26").Value = 80
$ws.Range("F27").Value = 445
$ws.Range("F28").Value = 667
$ws.Range("F29").Value = 2757
$ws.Range("F31").Value = 1840
$ws.Range("F34").Value = 837
$ws.Range("F38").Value = 90
$ws.Range("F39").Value = 938
$ws.Range("F40").Value = 31
$ws.Range("F43").Value = 94
$ws.Range("F44").Value = 533
$ws.Range("F45").Value = 599
$ws.Range("F46").Value = 904
$ws.Range("C48").Value = "杭州·理想乡动漫展-同人创作者大会"
$ws.Range("D48").Value = "阳城路雅澳杭州电商产业园西侧约200米 杭州大会展中心"
$ws.Range("E48").Value = "2024.09.15 10:00-09.16 17:00"
$ws.Range("F48").Value = 4239
$ws.Range("G48").Value = 1
$ws.Range("H48").Value = "https://show.bilibili.com/platform/detail.html?id=83822"
$ws.Range("I48").Value = "//i2.hdslb.com/bfs/openplatform/202404/GGEZUjGw1711959030111.png"
$ws.Range("B49").Value = "2024-10-05"
$ws.Range("C49").Value = "杭州·鸢飞鱼跃代号鸢only"
$ws.Range("D49").Value = "望江东路333号 杭州瑞莱克斯大酒店"
$ws.Range("E49").Value = "2024.10.05 09:30-10.05 17:00"
$ws.Range("F49").Value = 88
$ws.Range("G49").Value = 85
$ws.Range("H49").Value = "https://show.bilibili.com/platform/detail.html?id=88452"
$ws.Range("I49").Value = "//i2.hdslb.com/bfs/openplatform/202406/etOXBCrl1719678030944.jpeg"
